$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("F5").Value = "часы"
$ws.Range("G5").Value = "проценты"

$ws.Activate()
$ws.Range("I6").Select()
$excel.ActiveWindow.Zoom = 120
